# PROS-12595 GSKAU-SAND template update
#
# 1. Validation_List defined name now points at the 1st (and only) external
#    workbook reference instead of the stale "[2]" index.
# 2. The "Facings SOS" / "Availability" scene-type cell (C2/C3) text is
#    updated to list the new main-shelf scene types ahead of the existing
#    "- Grcy" variants, and gets comma separators between lines.
# 3. Row 2 / Row 3 grow taller (and a handful of column widths widen
#    slightly) to accommodate the longer wrapped text; the active
#    selection on the sheet moves from D13 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up the external-reference index used by Validation_List -------
$validationName = $wb.Names.Item("Validation_List")
$validationName.RefersTo = '=[1]Set_up!$A$90:$A$124'

# --- 2. Update the wrapped "Scene type / Tasks" text in C2 and C3 ---------
$newSceneText = "Pain Main Shelf, Oral Main Shelf, Respiratory Main Shelf, NRT Main Shelf, Other Main Shelf, Pain Main Shelf - Grcy,`nOral Main Shelf - Grcy,`nNRT Main Shelf - Grcy,`nCounter Unit - Grcy,`nFloor Bin - Grcy,`nGondola End - Grcy,`nHangsell - Grcy,`nClipstrip - Grcy,`nHotspot Tray – Grcy"

$ws.Range("C2").Value = $newSceneText
$ws.Range("C3").Value = $newSceneText

# --- 3. Resize rows/columns to match the re-wrapped content ---------------
# RowHeight takes points directly (no internal re-quantisation observed).
$ws.Rows.Item(2).RowHeight = 326.45
$ws.Rows.Item(3).RowHeight = 164.1

# ColumnWidth is stored internally in sixths of a character, so back out the
# 5/6-character padding before assigning to land as close as possible on the
# target stored width.
$ws.Columns.Item(1).ColumnWidth = 21.55533063427797    # A: 21.9595141700405 -> 22.3886639676113
$ws.Columns.Item(3).ColumnWidth = 30.016869095816467   # C: 30.4210526315789 -> 30.8502024291498
$ws.Columns.Item(7).ColumnWidth = 21.66059379217277    # G: 22.1740890688259 -> 22.4939271255061
$ws.Columns.Item(10).ColumnWidth = 21.87516869095817   # J: 22.4939271255061 -> 22.7085020242915
$ws.Columns.Item(11).ColumnWidth = 29.587719298245567  # K: 29.9919028340081 -> 30.4210526315789

# --- 4. Move the active selection from D13 to C3 ---------------------------
$null = $ws.Range("C3").Select()
